$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Proj2 (row 6) autograder result came in: 41/41.
$ws.Range("D6").Value = 41
$ws.Range("E6").Value = 41

# F6 = D6/E6, percentage-formatted like the other "%" cells above it (F3:F5)
$ws.Range("F6").Formula = "=D6/E6"
$ws.Range("F6").NumberFormat = "0.00%"

# H6 = G6*F6 (weighted value), newly filled in now that F6 has a real value
$ws.Range("H6").Formula = "=G6*F6"

# Proj2's weight (G6) now counts toward the "available points so far" total
$ws.Range("J7").Formula = "=SUM(G3:G6)"

# leave the selection on F6, matching the saved view
$ws.Range("F6").Select()
